$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns F:V) between rows 116 and 117 ---
$row116 = $ws.Range("F116:V116").Value2
$row117 = $ws.Range("F117:V117").Value2
$ws.Range("F116:V116").Value = $row117
$ws.Range("F117:V117").Value = $row116

# --- Swap match data (columns F:V) between rows 122 and 123 ---
$row122 = $ws.Range("F122:V122").Value2
$row123 = $ws.Range("F123:V123").Value2
$ws.Range("F122:V122").Value = $row123
$ws.Range("F123:V123").Value = $row122

# --- Append new row 125 (copy formatting from row 124, then set values) ---
$ws.Range("A124:V124").Copy($ws.Range("A125:V125"))

$ws.Range("A125").Value = 124
$ws.Range("B125").Value = "bulgaria"
$ws.Range("C125").Value = "vtora-liga"
$ws.Range("D125").Value = "2023-2024"
$ws.Range("E125").Value = 45224.57291666666
$ws.Range("F125").Value = "Litex Lovech"
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = "Spartak Varna"
$ws.Range("I125").Value = 2
$ws.Range("J125").Value = 3.45
$ws.Range("K125").Value = "24/10/2023 01:12"
$ws.Range("L125").Value = 4.95
$ws.Range("M125").Value = "25/10/2023 12:24"
$ws.Range("N125").Value = 2.84
$ws.Range("O125").Value = "24/10/2023 01:12"
$ws.Range("P125").Value = 3.57
$ws.Range("Q125").Value = "25/10/2023 12:24"
$ws.Range("R125").Value = 1.98
$ws.Range("S125").Value = "24/10/2023 01:12"
$ws.Range("T125").Value = 1.6
$ws.Range("U125").Value = "25/10/2023 11:25"
$ws.Range("V125").Value = "https://www.betexplorer.com/football/bulgaria/vtora-liga/litex-lovech-spartak-varna/Ox2inXmr/"
